$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the "Price" column to Text format so numeric-looking values
# (e.g. "1.000", "313.92") are stored verbatim as text, matching the
# source data which uses inline strings rather than numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.083.74"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "1.892.09"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "313.92"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.5053"
$ws.Range("E7").Value = "  -2.13%  "
$ws.Range("D8").Value = "0.3896"
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("D9").Value = "0.09249"
$ws.Range("E9").Value = "  -5.31%  "
$ws.Range("D10").Value = "1.128"
$ws.Range("E10").Value = "  -2.70%  "
$ws.Range("D11").Value = "41.84"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("E12").Value = "  -2.52%  "
$ws.Range("D13").Value = "20.81"
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("D14").Value = "1.897.61"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "7.291"
$ws.Range("E15").Value = "  -3.98%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "92.27"
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("D18").Value = "0.00001107"
$ws.Range("E18").Value = "  -3.08%  "
$ws.Range("D19").Value = "0.06654"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "17.85"
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "6.211"
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("D23").Value = "28.127.35"
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("D24").Value = "11.40"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").Value = "2.321"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("D26").Value = "2.107.61"
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("D27").Value = "2.544"
$ws.Range("E27").Value = "  -5.75%  "
$ws.Range("D28").Value = "158.54"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("D29").Value = "20.82"
$ws.Range("E29").Value = "  -2.14%  "
$ws.Range("D30").Value = "126.94"
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("D31").Value = "1.078"
$ws.Range("E31").Value = "  -2.77%  "
$ws.Range("E32").Value = "  -2.71%  "
$ws.Range("D33").Value = "5.608"
$ws.Range("E33").Value = "  -2.79%  "
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("D35").Value = "9.463"
$ws.Range("E35").Value = "  -4.19%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.06606"
$ws.Range("E36").Value = "  -3.02%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "1.343"
$ws.Range("E37").Value = "  +12.34%  "
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").Value = "0.2197"
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("E40").Value = "  -4.21%  "
$ws.Range("D41").Value = "11.67"
$ws.Range("E41").Value = "  -1.82%  "
$ws.Range("D42").Value = "0.6438"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").Value = "4.969"
$ws.Range("E43").Value = "  -3.17%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "13.30"
$ws.Range("E45").Value = "  -3.02%  "
$ws.Range("D46").Value = "0.6049"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("E47").Value = "  +1.57%  "
$ws.Range("D48").Value = "3.688"
$ws.Range("E48").Value = "  -2.50%  "
$ws.Range("D49").Value = "2.003"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").Value = "122.10"
$ws.Range("E50").Value = "  -2.56%  "
$ws.Range("D51").Value = "1.194"
$ws.Range("E51").Value = "  -1.84%  "

# Restore the default (no explicit) style on the Price column so the
# output matches the original workbook formatting.
$ws.Range("D2:D51").Style = "Normal"
